# Daily attendance processing - 2025-12-03 22:27:30
#
# For every row in the "Recorded By" column (G), whenever the value is a
# comma-separated list that contains "System" but does not already start
# with "System", move "System" to the front by reversing the order of the
# comma-separated entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value()

    if ($null -eq $value) { continue }

    $text = [string]$value
    if ($text -notmatch ",") { continue }

    $parts = $text -split "," | ForEach-Object { $_.Trim() }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p -eq "System") { $hasSystem = $true }
    }

    if (-not $hasSystem) { continue }
    if ($parts[0] -eq "System") { continue }

    $reversed = $parts[($parts.Count - 1)..0]
    $newText = [string]::Join(", ", $reversed)

    $cell.Value = $newText
}
